$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.070706844329834
$ws.Range("B1").Value = 3.986205816268921
$ws.Range("C1").Value = 5.695708274841309
$ws.Range("D1").Value = 1.601129651069641
$ws.Range("E1").Value = 0.9609283804893494
